$d = $word.ActiveDocument

# 1) Old company name placeholder -> current company name with optional old-name suffix.
$d.Content.Find.Execute(
    "{companyOldName2}", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "{companyName} {#hasCompanyOldName}[{companyOldName}]{/hasCompanyOldName}", 2)

# 2) Shareholder names on the certificate -> use the certificate-specific name fields.
$d.Content.Find.Execute(
    "{#hasShareholder_1}{shareholderName_1}{/hasShareholder_1}", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "{#hasShareholder_1}{shareholderNameCertificate_1}{/hasShareholder_1}", 2)

$d.Content.Find.Execute(
    "{#hasShareholder_2}{shareholderName_2}{/hasShareholder_2}", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "{#hasShareholder_2}{shareholderNameCertificate_2}{/hasShareholder_2}", 2)

$d.Content.Find.Execute(
    "{#hasShareholder_3}{shareholderName_3}{/hasShareholder_3}", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "{#hasShareholder_3}{shareholderNameCertificate_3}{/hasShareholder_3}", 2)
